$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new timeline entry as row 23 (Day 20)
$ws.Range("A23").Value = 20
$ws.Range("B23").Value = "22/3/2024"
$ws.Range("C23").Value = 6.5
$ws.Range("D23").Value = "Refactored search functionality, a lot of features added (category search, etc)"

# Match the formatting of the preceding data row (centered style, same as rows 4-22)
$ws.Range("A22:D22").Copy() | Out-Null
$ws.Range("A23:D23").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = 0

# Move the active cell selection from C25 to D25, as recorded in the sheet view
$ws.Range("D25").Select()
